# Adapted WebPage and stored auth info as base64-encoding in config
#
# The underlying OOXML diff for this commit shows the workbook being
# re-saved (version/namespace "noise" from a different Excel build) plus
# one concrete content edit on the "Mapping CSV2openEHR" sheet: the
# placeholder mapping that had been filled into row 7 (column B the
# FLAT-Path selection, column C a stray 0) is removed again, leaving row 7
# with only its CSV-Column label in A7 - matching every other still-unmapped
# row (rows 2-6, 8-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping CSV2openEHR")

# Clear out the accidentally-populated mapping (B7) and index (C7) cells.
$ws.Range("B7:C7").ClearContents()

# Reset the active selection back to the top of the sheet (A1), matching
# the saved cursor position in the target file (previously parked at D11).
$ws.Range("A1").Select() | Out-Null
